$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.509.02'
$ws.Range('E2').Value = '  -3.89%  '
$ws.Range('D3').Value = '2.464.87'
$ws.Range('E3').Value = '  -6.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.49'
$ws.Range('E5').Value = '  -4.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.98'
$ws.Range('E6').Value = '  -6.27%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.594'
$ws.Range('E8').Value = '  -5.54%  '
$ws.Range('D9').Value = '2.462.60'
$ws.Range('E9').Value = '  -6.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -9.37%  '
$ws.Range('E11').Value = '  -6.19%  '
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.355'
$ws.Range('E13').Value = '  -7.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.24'
$ws.Range('E14').Value = '  -8.64%  '
$ws.Range('D15').Value = '2.911.86'
$ws.Range('E15').Value = '  -6.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000167'
$ws.Range('E16').Value = '  -9.79%  '
$ws.Range('D17').Value = '61.425.87'
$ws.Range('E17').Value = '  -3.92%  '
$ws.Range('D18').Value = '2.467.26'
$ws.Range('E18').Value = '  -6.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.12'
$ws.Range('E19').Value = '  -8.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.18'
$ws.Range('E20').Value = '  -7.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.20'
$ws.Range('E21').Value = '  -7.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '319.99'
$ws.Range('E22').Value = '  -7.32%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.88'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.34'
$ws.Range('E25').Value = '  -5.79%  '
$ws.Range('D26').Value = '0.0₃0992'
$ws.Range('E26').Value = '  -11.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '564.00'
$ws.Range('E27').Value = '  -3.13%  '
$ws.Range('D28').Value = '2.606.39'
$ws.Range('E28').Value = '  -6.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.50'
$ws.Range('E29').Value = '  -8.71%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.30'
$ws.Range('E31').Value = '  -10.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.69'
$ws.Range('E32').Value = '  -6.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.149'
$ws.Range('E33').Value = '  -7.84%  '
$ws.Range('E34').Value = '  -6.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  -7.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.92'
$ws.Range('E36').Value = '  -10.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.87'
$ws.Range('E37').Value = '  -10.89%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.380'
$ws.Range('E39').Value = '  -5.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.46'
$ws.Range('E40').Value = '  -6.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.46'
$ws.Range('E41').Value = '  -6.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.75'
$ws.Range('E42').Value = '  -8.57%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.58'
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  -5.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '146.93'
$ws.Range('E46').Value = '  -9.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.62'
$ws.Range('E47').Value = '  -7.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.66'
$ws.Range('E48').Value = '  -10.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0538'
$ws.Range('E49').Value = '  -8.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.594'
$ws.Range('E50').Value = '  -6.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0941'
$ws.Range('E51').Value = '  -6.24%  '
